$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 391.95456
$ws.Range("I33").Value = 419
$ws.Range("K33").Value = 419
$ws.Range("M33").Value = -190

$ws.Range("H97").Value = 1838.25
$ws.Range("J97").Value = 1838.25
$ws.Range("L97").Value = 5514.75
$ws.Range("N97").Value = -6506.75

$ws.Range("H103").Value = 2676.889
$ws.Range("I103").Value = 749.25
$ws.Range("J103").Value = 4219
$ws.Range("K103").Value = 2247.75
$ws.Range("L103").Value = 12657
$ws.Range("M103").Value = -1661.75
$ws.Range("N103").Value = -13829

$ws.Range("H132").Value = 792.3333
$ws.Range("I132").Value = 686.5135
$ws.Range("K132").Value = 2059.5405
$ws.Range("M132").Value = 470.4594999999999

$ws.Range("H137").Value = 13516252
$ws.Range("I137").Value = 52633372
$ws.Range("J137").Value = 3065.2
$ws.Range("K137").Value = 157900116
$ws.Range("L137").Value = 9195.599999999999
$ws.Range("M137").Value = -157897566
$ws.Range("N137").Value = -14295.6

$ws.Range("H138").Value = 2844.247
$ws.Range("I138").Value = 1652.409
$ws.Range("J138").Value = 3235.597
$ws.Range("K138").Value = 4957.227000000001
$ws.Range("L138").Value = 9706.791000000001
$ws.Range("M138").Value = 182.7729999999992
$ws.Range("N138").Value = -19986.791

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3925.4658
$ws.Range("I32").Value = 2768.175
$ws.Range("K32").Value = 2768.175
$ws.Range("M32").Value = -2481.175

$ws.Range("H45").Value = 100001750
$ws.Range("I45").Value = 100001750
$ws.Range("K45").Value = 100001750
$ws.Range("M45").Value = -100001373

$ws.Range("H61").Value = 27007
$ws.Range("J61").Value = 29014
$ws.Range("L61").Value = 29014
$ws.Range("N61").Value = -29438

$ws.Range("H110").Value = 6269.125
$ws.Range("I110").Value = 5607.684
$ws.Range("K110").Value = 5607.684
$ws.Range("M110").Value = -3562.684

$ws.Range("H136").Value = 27007
$ws.Range("J136").Value = 29014
$ws.Range("L136").Value = 87042
$ws.Range("N136").Value = -92142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4534.1035
$ws.Range("J20").Value = 5387.4443
$ws.Range("L20").Value = 5387.4443
$ws.Range("N20").Value = -5881.4443

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H134").Value = 3646.3
$ws.Range("I134").Value = 2181.125
$ws.Range("K134").Value = 6543.375
$ws.Range("M134").Value = -4008.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 8996
$ws.Range("I22").Value = 9986.5
$ws.Range("J22").Value = 8599.799999999999
$ws.Range("K22").Value = 9986.5
$ws.Range("L22").Value = 8599.799999999999
$ws.Range("M22").Value = -9636.5
$ws.Range("N22").Value = -9299.799999999999

$ws.Range("H31").Value = 27678.404
$ws.Range("I31").Value = 2119.6667
$ws.Range("K31").Value = 2119.6667
$ws.Range("M31").Value = -1824.6667

$ws.Range("H34").Value = 27678.404
$ws.Range("I34").Value = 2119.6667
$ws.Range("K34").Value = 2119.6667
$ws.Range("M34").Value = -1917.6667

$ws.Range("H58").Value = 3942.9333
$ws.Range("I58").Value = 1979.8334
$ws.Range("K58").Value = 1979.8334
$ws.Range("M58").Value = -1776.8334

$ws.Range("H99").Value = 2874.5
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 2874.5
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 2874.5
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -5870.5

$ws.Range("H112").Value = 60789.25
$ws.Range("J112").Value = 60789.25
$ws.Range("L112").Value = 60789.25
$ws.Range("N112").Value = -63743.25

$ws.Range("H126").Value = 2874.5
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 2874.5
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 8623.5
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -13563.5

$ws.Range("H132").Value = 3673.111
$ws.Range("I132").Value = 3155.8235
$ws.Range("K132").Value = 9467.470499999999
$ws.Range("M132").Value = -6937.470499999999

$ws.Range("H136").Value = 3942.9333
$ws.Range("I136").Value = 1979.8334
$ws.Range("K136").Value = 5939.5002
$ws.Range("M136").Value = -3389.5002

$ws.Range("H140").Value = 87638.30499999999
$ws.Range("J140").Value = 87638.30499999999
$ws.Range("L140").Value = 87638.30499999999
$ws.Range("N140").Value = -97998.30499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 6242178.5
$ws.Range("J122").Value = 7144668
$ws.Range("L122").Value = 64302012
$ws.Range("N122").Value = -64306912

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 458917.53
$ws.Range("I80").Value = 502908.7
$ws.Range("J80").Value = 19006
$ws.Range("K80").Value = 502908.7
$ws.Range("L80").Value = 19006
$ws.Range("M80").Value = -501910.7
$ws.Range("N80").Value = -21002

$ws.Range("H83").Value = 458917.53
$ws.Range("I83").Value = 502908.7
$ws.Range("J83").Value = 19006
$ws.Range("K83").Value = 2514543.5
$ws.Range("L83").Value = 95030
$ws.Range("M83").Value = -2509551.5
$ws.Range("N83").Value = -105014

$ws.Range("H102").Value = 2045.6072
$ws.Range("I102").Value = 1189.6
$ws.Range("K102").Value = 1189.6
$ws.Range("M102").Value = 432.4000000000001

$ws.Range("H126").Value = 3380.0908
$ws.Range("I126").Value = 1646
$ws.Range("K126").Value = 4938
$ws.Range("M126").Value = -2468

$ws.Range("H132").Value = 2599.524
$ws.Range("I132").Value = 2105.0303
$ws.Range("K132").Value = 6315.090899999999
$ws.Range("M132").Value = -3785.090899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7781.4136
$ws.Range("I40").Value = 6711.9443
$ws.Range("K40").Value = 6711.9443
$ws.Range("M40").Value = -6575.9443

$ws.Range("H82").Value = 5685.077
$ws.Range("I82").Value = 1521
$ws.Range("J82").Value = 10543.167
$ws.Range("K82").Value = 1521
$ws.Range("L82").Value = 10543.167
$ws.Range("M82").Value = -1160
$ws.Range("N82").Value = -11265.167

$ws.Range("H85").Value = 5685.077
$ws.Range("I85").Value = 1521
$ws.Range("J85").Value = 10543.167
$ws.Range("K85").Value = 1521
$ws.Range("L85").Value = 10543.167
$ws.Range("M85").Value = -273
$ws.Range("N85").Value = -13039.167

$ws.Range("H122").Value = 271605.2
$ws.Range("I122").Value = 336425.1
$ws.Range("K122").Value = 1009275.3
$ws.Range("M122").Value = -1006825.3

$ws.Range("H132").Value = 5709.2827
$ws.Range("I132").Value = 6027.353
$ws.Range("J132").Value = 4808.0835
$ws.Range("K132").Value = 18082.059
$ws.Range("L132").Value = 14424.2505
$ws.Range("M132").Value = -15552.059
$ws.Range("N132").Value = -19484.2505

$ws.Range("H136").Value = 2843.426
$ws.Range("I136").Value = 1752.3077
$ws.Range("K136").Value = 5256.9231
$ws.Range("M136").Value = -2706.9231

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 999.8
$ws.Range("I100").Value = 813.2727
$ws.Range("J100").Value = 1512.75
$ws.Range("K100").Value = 1626.5454
$ws.Range("L100").Value = 3025.5
$ws.Range("M100").Value = -1085.5454
$ws.Range("N100").Value = -4107.5

$ws.Range("H122").Value = 1979.093
$ws.Range("I122").Value = 1431.1562
$ws.Range("K122").Value = 4293.4686
$ws.Range("M122").Value = -1843.4686

$ws.Range("H132").Value = 2394.6667
$ws.Range("I132").Value = 1755.1923
$ws.Range("K132").Value = 5265.5769
$ws.Range("M132").Value = -2735.5769
